# Auto-generated update of cryptos list values (prices & volume%) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.357.35'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.879.69'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = '''242.32'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '''0.08041'
$ws.Range('E8').Value = '  +4.02%  '
$ws.Range('D9').Value = '''0.3128'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('D10').Value = '''25.28'
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('D11').Value = '''0.08363'
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('D12').Value = '1.876.77'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '''5.258'
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').Value = '''91.43'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').Value = '''6.279'
$ws.Range('E16').Value = '  +5.01%  '
$ws.Range('D17').Value = '''0.000008383'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '29.364.81'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = '''241.00'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').Value = '2.134.71'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D23').Value = '''7.810'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '''0.1589'
$ws.Range('E25').Value = '  -1.96%  '
$ws.Range('D26').Value = '''163.11'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '''9.068'
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').Value = '''1.508'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '''4.423'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D32').Value = '''1.200'
$ws.Range('E32').Value = '  -6.04%  '
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('D34').Value = '''1.951'
$ws.Range('E34').Value = '  +1.57%  '
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').Value = '''0.7510'
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').Value = '''2.699'
$ws.Range('E37').Value = '  +0.64%  '
$ws.Range('D38').Value = '1.294.94'
$ws.Range('E38').Value = '  +11.53%  '
$ws.Range('D39').Value = '''0.01887'
$ws.Range('E39').Value = '  +1.50%  '
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').Value = '''6.586'
$ws.Range('E41').Value = '  +3.57%  '
$ws.Range('D42').Value = '''110.68'
$ws.Range('E42').Value = '  +3.88%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''0.8928'
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''73.09'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('E45').Value = '  +8.38%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '2.023.37'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').Value = '''0.5205'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').Value = '''9.485'
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').Value = '''0.4366'
$ws.Range('E51').Value = '  +1.60%  '
